# Ticket 46 - Add "Answer:" / "${answerToLifeTheUniverseAndEverything}" and
# "Pick A Card:" / "${jett:pickACard()}" rows beneath the existing "EscExprs:" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Answer:"
$ws.Range("B31").Value = '${answerToLifeTheUniverseAndEverything}'

$ws.Range("A32").Value = "Pick A Card:"
$ws.Range("B32").Value = '${jett:pickACard()}'
